# "Scheduled" posts in rows 21-38 were originally queued with a placeholder
# send time (45337.004861111112 == 2024-02-15 00:07). The new APScheduler
# job-completion listener (my_listener) confirms/logs the real dispatch
# time for these jobs, so update the Scheduled Time column to the actual
# time the jobs fired (45337.501388888886 == 2024-02-15 12:02) and leave
# the worksheet selection parked on the range that was just refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newScheduledTime = 45337.501388888886

$ws.Range("F21:F38").Value2 = $newScheduledTime

$ws.Range("F21:F38").Select()
